# SAM TODO.xlsx - capital-costs / self-shading bugfix entry
#
# "PV self-shading testing" (row 34) is now Done, and a new follow-up
# task is inserted right after it: "Re-arrange self-shading inputs in UI
# with system design? Check inputs for usability in SDK" (owner Janine,
# status Not done). Everything that used to be below row 34 shifts down
# by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 35 - native Excel row-insert semantics: formats
# copy down from the row above, formulas/ranges below get re-pointed
# automatically, and existing rows 35..61 become 36..62.
$ws.Rows("35:35").Insert() | Out-Null

# Row 34 ("PV self-shading testing") flips from "Not done" to "Done".
$ws.Cells.Item(34, 1).Value2 = "Done"

# Populate the newly inserted row 35 with the follow-up task.
$ws.Cells.Item(35, 1).Value2 = "Not done"
$ws.Cells.Item(35, 2).Value2 = "Re-arrange self-shading inputs in UI with system design? Check inputs for usability in SDK"
$ws.Cells.Item(35, 3).Value2 = "Janine"
# D35 (hours) and E35 (priority) stay blank for this new entry.

# The running total in H17 summed D17:D39; with the extra row it now
# needs to reach down to D40.
$ws.Range("H17").Formula = "=SUM(D17:D40)"

# Restore the view: scrolled down so row 28 is at the top, cursor on B36.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("B36").Select() | Out-Null
